$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchResults")

$ws.Range("A4").Value = "search_3"
$ws.Range("A5").Value = "search_4"
$ws.Range("B4").Value = "pants"
$ws.Range("B5").Value = "TAURUS"
$ws.Range("B3").Value = "Backpack"

$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 1

$ws.Range("C3").Select()
